$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "CRUD Rendez-vous" task as done
$ws.Range("C25").Value = "YES"

# Update the secretary home-page task description
$ws.Range("A26").Value = 'Accueil "ActiveVisitsForPets" et "VeterinaryAvailabilityForWeek"'

# Update the active selection on the sheet
$ws.Range("B11").Select()
